# New profile script implementation
# Adds a new "Profile70" test case row (row 71) to the Profile worksheet,
# mirroring the formatting of the row above it (row 70), and updates the
# sheet's active selection to reflect where the author ended up working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the formatting of the last existing data row (row 70) down onto the
# new row (row 71) so the new row's cell styles match (s="2"/"7"/"11"/"2"/"7").
$ws.Range("A70:E70").Copy() | Out-Null
$ws.Range("A71:E71").PasteSpecial(-4122) | Out-Null

# Populate the new test case values.
$ws.Range("A71").Value = "Profile70"
$ws.Range("B71").Value = "OPQA-2109||OPQA-2110"
$ws.Range("C71").Value = "Verify that user has the ability to add or change their profile photo from the profile modal.||Verify that user has the ability to delete his profile photo from the profile modal."
$ws.Range("D71").Value = "Y"

# Row 71 wraps to two lines like the other multi-line rows, so give it the
# same row height (30) as its neighbours.
$ws.Rows.Item(71).RowHeight = 30

# Reflect the author's final cursor position/selection in the sheet view.
$ws.Range("J68").Select() | Out-Null
